# issue #5: stock data from json to db
# Adds "category", "source_file" and "index" columns to the 股票 (stock) sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")   # stock sheet

# --- Insert a new column I ("category") right after "property_category" (H) ---
# This shifts old I (date), J (legislator_name), K (legislator_id) to J, K, L.
# Insert() carries the neighbouring column's formatting along, so the new
# column I already ends up styled like the rest of the header/data cells.
$ws.Columns.Item(9).Insert()

# --- Header row (row 1) ---
$ws.Range("I1").Value = "category"

# Give the two brand-new trailing header cells (M1, N1) the same look as the
# existing header cells (bold + border), by copying the format from L1 first.
$ws.Range("L1").Copy()
$ws.Range("M1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("L1").Copy()
$ws.Range("N1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Data rows ---
# Row 2
$ws.Range("I2").Value = "normal"
$ws.Range("M2").Value = "tmp1fff1"
$ws.Range("N2").Value = 71

# Row 3
$ws.Range("I3").Value = "normal"
$ws.Range("M3").Value = "tmp1fff1"
$ws.Range("N3").Value = 72

# Row 4
$ws.Range("I4").Value = "normal"
$ws.Range("M4").Value = "tmp1fff1"
$ws.Range("N4").Value = 73
